# feat add api qr code + fix l'interface messagerie
#
# Appends 11 new question/answer rows (44-54) to the "Feuil1" chatbot
# dataset sheet. Cell values are written in the exact order required to
# reproduce the target shared-string table ordering (the author filled in
# column A for rows 48-50 before going back to fill B48, which is why B48's
# shared-string index is higher than the ones that follow it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44
$ws.Cells.Item(44, 1).Value = "iheb"
$ws.Cells.Item(44, 2).Value = "c'est un prof tres intelligent dans le domaine IT et il est l'encadreur des autres équipes"

# Row 45
$ws.Cells.Item(45, 1).Value = "kifeh nestaamel baladity"
$ws.Cells.Item(45, 2).Value = "baladity esta3melha sehel barcha jareb les fonctionnalités eli andna w ija koli "

# Row 46
$ws.Cells.Item(46, 1).Value = "fama des evenements ?"
$ws.Cells.Item(46, 2).Value = "ey tnajem tchouf les evenements eli mawjoudin l periode hedhy sehla yesser odkhel lel evenement hawka fel menu"

# Row 47
$ws.Cells.Item(47, 1).Value = "kifeh n3ady reclamation"
$ws.Cells.Item(47, 2).Value = "sehel barcha bech t3adi reclamation juste todkhel lel reclamation w takhtar type reclamation eli bech taadih w testana directeur yjewbouk al reclamation mte3ek"

# Row 48 (column A only for now; B48 is filled in further below)
$ws.Cells.Item(48, 1).Value = "plus mta3 l application hedhy chneya"

# Row 49
$ws.Cells.Item(49, 1).Value = "fama flous "
$ws.Cells.Item(49, 2).Value = "le betbi3a caisse fergha"

# Row 50
$ws.Cells.Item(50, 1).Value = "chnwa famma jdid"
$ws.Cells.Item(50, 2).Value = "bech taaref chnwa famma jdid odkhel lel actualite talka kol chy sayer "

# Now go back and fill in B48
$ws.Cells.Item(48, 2).Value = "tsahel alik l contact f kol baladeya"

# Row 51
$ws.Cells.Item(51, 1).Value = "kifeh nhabet publicite"
$ws.Cells.Item(51, 2).Value = "easy peasy taamer formlaire w tet3ada lel paiement wtkhales b soum ramzi w jawek foll"

# Row 52
$ws.Cells.Item(52, 1).Value = "kadech tahky men lougha"
$ws.Cells.Item(52, 2).Value = "je parle en francais,arabe,anglais"

# Row 53
$ws.Cells.Item(53, 1).Value = "chneya loumour"
$ws.Cells.Item(53, 2).Value = "famm chyy jdid"

# Row 54
$ws.Cells.Item(54, 1).Value = "kadeh omrek"
$ws.Cells.Item(54, 2).Value = "En tant que chatbot maandich reponse lel question mte3ek"

# Leave the selection on the last cell entered, matching the author's
# final cursor position in the saved workbook.
$ws.Range("B54").Select()
